# ============================================================================
# Edit script: applies the commit "fix lỗi trong report cơ sở. Thêm cột ghi
# chú trong báo cáo về chi tiêu" to the workbook.
#
#   1. Rename sheet "Thưởng" -> "Đơn thu nợ" and replace its contents with a
#      debt-collection ("thu nợ") order report (25 columns).
#   2. Sheet "Đơn phụ phẫu 1": add 3 new service rows and refresh the total
#      row.
#   3. Sheet "Lương": insert a "Chiết khấu thu nợ tại <cơ sở>" row for each
#      location, drop the old "Thưởng tại CẦN THƠ" row, and refresh every
#      downstream total.
# ============================================================================

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    if ($null -eq $text) {
        $cell.Value = ""
    } else {
        $cell.Value = $text
    }
}

function Set-NumCell($ws, $row, $col, $num) {
    $ws.Cells.Item($row, $col).Value = $num
}

# ----------------------------------------------------------------------------
# 1) Sheet "Đơn phụ phẫu 1" (sheet #2) - add three new rows before the total
# ----------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Đơn phụ phẫu 1")

# push the existing "Tổng" row (row 5) down to row 8, three times
$ws2.Rows.Item(5).Insert()
$ws2.Rows.Item(5).Insert()
$ws2.Rows.Item(5).Insert()

# row 5
Set-TextCell $ws2 5 1 "HD-LUXURY"
Set-NumCell  $ws2 5 2 630
Set-TextCell $ws2 5 3 "08-05-2024"
Set-TextCell $ws2 5 4 "CẦN THƠ"
Set-TextCell $ws2 5 5 $null
Set-TextCell $ws2 5 6 "CTV"
Set-TextCell $ws2 5 7 "Nâng mũi"
Set-TextCell $ws2 5 8 "Lâm Hoàng Phú"
Set-NumCell  $ws2 5 9 100000

# row 6
Set-TextCell $ws2 6 1 "HD-LUXURY"
Set-NumCell  $ws2 6 2 633
Set-TextCell $ws2 6 3 "08-07-2024"
Set-TextCell $ws2 6 4 "CẦN THƠ"
Set-TextCell $ws2 6 5 "Bạch Nhi"
Set-TextCell $ws2 6 6 "Cá nhân"
Set-TextCell $ws2 6 7 "Nâng mũi"
Set-TextCell $ws2 6 8 "Lâm Hoàng Phú"
Set-NumCell  $ws2 6 9 100000

# row 7
Set-TextCell $ws2 7 1 "HD-LUXURY"
Set-NumCell  $ws2 7 2 635
Set-TextCell $ws2 7 3 "08-08-2024"
Set-TextCell $ws2 7 4 "CẦN THƠ"
Set-TextCell $ws2 7 5 "Nguyễn Bích Thuỳ"
Set-TextCell $ws2 7 6 "Khách cũ giới thiệu"
Set-TextCell $ws2 7 7 "Cắt mí"
Set-TextCell $ws2 7 8 "Lâm Hoàng Phú"
Set-NumCell  $ws2 7 9 50000

# row 8 ("Tổng") - refresh counts/totals (6 orders, 500000 total)
Set-TextCell $ws2 8 1 "Tổng"
Set-NumCell  $ws2 8 2 6
Set-NumCell  $ws2 8 9 500000

# ----------------------------------------------------------------------------
# 2) Sheet "Thưởng" -> "Đơn thu nợ" (sheet #3) - wholesale content swap
# ----------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Thưởng")
$ws3.Cells.Clear()
$ws3.Name = "Đơn thu nợ"

Set-TextCell $ws3 1 1  "Tiền tố"
Set-TextCell $ws3 1 2  "Mã đơn thu nợ"
Set-TextCell $ws3 1 3  "Lượng thu"
Set-TextCell $ws3 1 4  "Ngày thu"
Set-TextCell $ws3 1 5  "Cơ sở"
Set-TextCell $ws3 1 6  "Đơn nợ"
Set-TextCell $ws3 1 7  "Tên dịch vụ"
Set-TextCell $ws3 1 8  "Khách hàng"
Set-TextCell $ws3 1 9  "Nguồn khách"
Set-TextCell $ws3 1 10 "Sale chính"
Set-TextCell $ws3 1 11 "Đơn giá gốc"
Set-TextCell $ws3 1 12 "Sale phụ"
Set-TextCell $ws3 1 13 "Upsale"
Set-TextCell $ws3 1 14 "Đơn giá"
Set-TextCell $ws3 1 15 "Đã thanh toán"
Set-TextCell $ws3 1 16 "Bác sĩ 1"
Set-TextCell $ws3 1 17 "Bác sĩ 2"
Set-TextCell $ws3 1 18 "Tỉ lệ chiết khấu sale chính"
Set-TextCell $ws3 1 19 "Chiết khấu sale chính"
Set-TextCell $ws3 1 20 "Tỉ lệ chiết khấu sale phụ"
Set-TextCell $ws3 1 21 "Chiết khấu sale phụ"
Set-TextCell $ws3 1 22 "Tỉ lệ chiết khấu bác sĩ 1"
Set-TextCell $ws3 1 23 "Chiết khấu bác sĩ 1"
Set-TextCell $ws3 1 24 "Tỉ lệ chiết khấu bác sĩ 2"
Set-TextCell $ws3 1 25 "Chiết khấu bác sĩ 2"

# row 2 - the single debt-collection order
Set-TextCell $ws3 2 1  "TN"
Set-NumCell  $ws3 2 2  181
Set-NumCell  $ws3 2 3  1500000
Set-TextCell $ws3 2 4  "08-09-2024"
Set-TextCell $ws3 2 5  "CẦN THƠ"
Set-TextCell $ws3 2 6  "HD-LUXURY-538"
Set-TextCell $ws3 2 7  "Nâng mũi"
Set-TextCell $ws3 2 8  "Ngô Xuân Nhi"
Set-TextCell $ws3 2 9  "Cá nhân"
Set-TextCell $ws3 2 10 "Lâm Hoàng Phú"
Set-NumCell  $ws3 2 11 10000000
Set-TextCell $ws3 2 12 "Đỗ Thị Huyền Trân"
Set-NumCell  $ws3 2 13 8000000
Set-NumCell  $ws3 2 14 18000000
Set-NumCell  $ws3 2 15 11000000
Set-TextCell $ws3 2 16 "Lâm Thị Mỹ Hằng"
Set-TextCell $ws3 2 17 $null
Set-NumCell  $ws3 2 18 0.1
Set-NumCell  $ws3 2 19 110000
Set-NumCell  $ws3 2 20 0
Set-NumCell  $ws3 2 21 0
Set-NumCell  $ws3 2 22 0
Set-NumCell  $ws3 2 23 0
Set-NumCell  $ws3 2 24 0
Set-NumCell  $ws3 2 25 0

# row 3 - "Tổng"
Set-TextCell $ws3 3 1  "Tổng"
Set-NumCell  $ws3 3 2  1
Set-NumCell  $ws3 3 3  1500000
Set-NumCell  $ws3 3 11 10000000
Set-NumCell  $ws3 3 13 8000000
Set-NumCell  $ws3 3 14 18000000
Set-NumCell  $ws3 3 15 11000000
Set-NumCell  $ws3 3 18 0
Set-NumCell  $ws3 3 19 110000
Set-NumCell  $ws3 3 20 0
Set-NumCell  $ws3 3 21 0
Set-NumCell  $ws3 3 22 0
Set-NumCell  $ws3 3 23 0
Set-NumCell  $ws3 3 24 0
Set-NumCell  $ws3 3 25 0

# ----------------------------------------------------------------------------
# 3) Sheet "Lương" (sheet #4) - insert "Chiết khấu thu nợ tại <cơ sở>" rows
# ----------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Lương")

# Insert a fresh row right before each "Ứng lương tại <cơ sở>" row. Doing
# this from the bottom-most insertion point up keeps the earlier row numbers
# stable while we work.
$ws4.Rows.Item(31).Insert()   # before "Ứng lương tại SÓC TRĂNG" (was row 31)
$ws4.Rows.Item(21).Insert()   # before "Ứng lương tại LONG XUYÊN" (was row 21)
$ws4.Rows.Item(10).Insert()   # before "Ứng lương tại CẦN THƠ" (was row 10)

# CẦN THƠ block
Set-NumCell  $ws4 2 2 7
Set-NumCell  $ws4 3 2 812500
Set-NumCell  $ws4 8 2 500000
Set-TextCell $ws4 10 1 "Chiết khấu thu nợ tại CẦN THƠ"
Set-NumCell  $ws4 10 2 110000
Set-TextCell $ws4 11 1 "Ứng lương tại CẦN THƠ"
Set-NumCell  $ws4 11 2 -300000

# row 12 used to be "Thưởng tại CẦN THƠ" (4000000) - it no longer exists;
# row 12 is now "Tổng công tại LONG XUYÊN", already correct (value 0).

# LONG XUYÊN block
Set-TextCell $ws4 21 1 "Chiết khấu thu nợ tại LONG XUYÊN"
Set-NumCell  $ws4 21 2 0
Set-TextCell $ws4 22 1 "Ứng lương tại LONG XUYÊN"
Set-NumCell  $ws4 22 2 0

# SÓC TRĂNG block
Set-TextCell $ws4 32 1 "Chiết khấu thu nợ tại SÓC TRĂNG"
Set-NumCell  $ws4 32 2 0
Set-TextCell $ws4 33 1 "Ứng lương tại SÓC TRĂNG"
Set-NumCell  $ws4 33 2 0

# Totals
Set-TextCell $ws4 34 1 "Tổng lương tại CẦN THƠ"
Set-NumCell  $ws4 34 2 1122500
Set-TextCell $ws4 35 1 "Tổng lương tại LONG XUYÊN"
Set-NumCell  $ws4 35 2 0
Set-TextCell $ws4 36 1 "Tổng lương tại SÓC TRĂNG"
Set-NumCell  $ws4 36 2 0
Set-TextCell $ws4 37 1 "Tổng lương tại HỆ THỐNG"
Set-NumCell  $ws4 37 2 1122500
